$wb = $excel.ActiveWorkbook

# --- Sheet1: add the two new rows (A9/B9 numeric, A10/B10 the new "asd" string) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A9").Value = 123
$ws1.Range("B9").Value = 123
$ws1.Range("A10").Value = "asd"
$ws1.Range("B10").Value = "asd"

# Sheet1's view settles on A9 at 205% zoom (no longer the selected tab)
$ws1.Range("A9").Select()
$excel.ActiveWindow.Zoom = 205

# --- Sheet2: add C6 ("asd"); this sheet becomes the active / tab-selected one ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("C6").Value = "asd"
$ws2.Activate()
$ws2.Range("C6").Select()
$excel.ActiveWindow.Zoom = 160
